$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing segment names (currently in column A, rows 2-20) ---
$names = @()
for ($r = 2; $r -le 20; $r++) {
    $names += $ws.Cells.Item($r, 1).Text
}

# --- Insert a new column before column B, shifting the old PercActivations / ---
# --- PercSegmentAreas columns (old B, C) to the right (new C, D).           ---
$ws.Columns.Item(2).Insert(-4161)  # -4161 = xlShiftToRight

# --- The inserted column inherits column A's header-row formatting; the   ---
# --- new segment-name cells (B2:B20) should be plain/unstyled like the    ---
# --- old numeric data columns were.                                       ---
$ws.Range("B2:B20").Style = "Normal"

# --- Give the new header cell (B1) the same "header" style used by the     ---
# --- other header cells (bold, bordered, centered), then set its text.     ---
$ws.Cells.Item(1, 2).Font.Bold = $true
$ws.Cells.Item(1, 2).HorizontalAlignment = -4108
$ws.Cells.Item(1, 2).VerticalAlignment = -4160
$ws.Cells.Item(1, 2).Borders.LineStyle = 1
$ws.Cells.Item(1, 2).Value = "segments"

# --- Column A becomes the 0-based numeric segment index; column B gets the ---
# --- segment name text that used to live in column A.                     ---
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
}
